# Duplicate the six original class sheets, appending "1" to a first pass
# over all six, then appending "2" to a second pass over the first two
# (呂彥臻 and 王大同) — matching the new tabs added in the target workbook:
#   呂彥臻1, 王大同1, 楊小明1, 劉大象1, 陳中一1, 孫二1, 呂彥臻2, 王大同2

$wb = $excel.ActiveWorkbook

# Source sheet index (1-based, into the ORIGINAL 6 sheets) and the suffix
# to append to both the tab name and the text in A1, in the exact order
# the new sheets must appear.
$plan = @(
    @{ Src = 1; Suffix = "1" },
    @{ Src = 2; Suffix = "1" },
    @{ Src = 3; Suffix = "1" },
    @{ Src = 4; Suffix = "1" },
    @{ Src = 5; Suffix = "1" },
    @{ Src = 6; Suffix = "1" },
    @{ Src = 1; Suffix = "2" },
    @{ Src = 2; Suffix = "2" }
)

foreach ($item in $plan) {
    $srcSheet = $wb.Worksheets.Item($item.Src)
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

    # Copy the source sheet to the very end of the workbook.
    $srcSheet.Copy($null, $lastSheet)
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

    $baseName = $srcSheet.Name
    $baseTitle = $srcSheet.Range("A1").Text

    $newSheet.Name = $baseName + $item.Suffix
    $newSheet.Range("A1").Value = $baseTitle + $item.Suffix
}

# Restore the originally active sheet (the repeated Copy() calls above
# leave the newest sheet active/selected).
$wb.Worksheets.Item(1).Activate()

